$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Change A: shared string "2016-10-21 04:53:42" -> "2016-10-21 04:55:30"
#   Used by Overview!G2:G3 (Latest HO Xliff Generate Date)
#        and de-de!H2:H3   (Correspond Handoff Datetime)
$wsOverview.Range("G2").Value = "2016-10-21 04:55:30"
$wsOverview.Range("G3").Value = "2016-10-21 04:55:30"
$wsDeDe.Range("H2").Value = "2016-10-21 04:55:30"
$wsDeDe.Range("H3").Value = "2016-10-21 04:55:30"

# Change B: shared string "ht" -> "mt"
#   Used by zh-cn!E2:E3 AND de-de!E2:E3 (Priority)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# Change C: shared string "2016-10-21 04:53:10" -> "2016-10-21 04:55:18"
#   Used by zh-cn!H2:H3 (Correspond Handoff Datetime)
$wsZhCn.Range("H2").Value = "2016-10-21 04:55:18"
$wsZhCn.Range("H3").Value = "2016-10-21 04:55:18"

# Change D: shared string "2016-10-21 04:54:20" -> "2016-10-21 04:56:00"
#   Used by zh-cn!K2:K3 (Correspond Handback DateTime)
$wsZhCn.Range("K2").Value = "2016-10-21 04:56:00"
$wsZhCn.Range("K3").Value = "2016-10-21 04:56:00"

# Change E: shared string "2016-10-21 04:54:38" -> "2016-10-21 04:56:19"
#   Used by de-de!K2:K3 (Correspond Handback DateTime)
$wsDeDe.Range("K2").Value = "2016-10-21 04:56:19"
$wsDeDe.Range("K3").Value = "2016-10-21 04:56:19"
